$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.908.10"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.425.78"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.25"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.93"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +18.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.449.94"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.843.43"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.447.76"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.68"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.437.01"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.70"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.84"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.87"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.91"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.61"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.521.52"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0780"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.26"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.50"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.14"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.855"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +9.61%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.19"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.602"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0535"
$ws.Range("E45").Value = "  -3.99%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "262.20"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.74"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.19"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.46"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.862.20"
$ws.Range("E51").Value = "  -1.70%  "
